$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Rows 3-6: B column "yes" -> "added to db" (status updated for completed runs)
$ws.Range("B3").Value = "added to db"
$ws.Range("B4").Value = "added to db"
$ws.Range("B5").Value = "added to db"
$ws.Range("B6").Value = "added to db"

# Row 7: "yes"/"running" -> "added to db"/"error"
$ws.Range("B7").Value = "added to db"
$ws.Range("C7").Value = "error"

# Row 8: "not started"/"not started" -> "added to db"/"error"
$ws.Range("B8").Value = "added to db"
$ws.Range("C8").Value = "error"

# Row 9: "not started"/"not started" -> "yes"/"running"
$ws.Range("B9").Value = "yes"
$ws.Range("C9").Value = "running"

# Update the active selection to reflect the last cell touched during the session
$ws.Range("B10").Select()

$wb.Save()
